# Automatic update of files.
# Remove the stray second argument ("test") from the HYPERLINK() formulas
# in row 2, columns S:Y, leaving just the single URL argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/artfynd/A 32298-2023.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/kartor/A 32298-2023.png")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/knärot/A 32298-2023.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/klagomål/A 32298-2023.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/klagomålsmail/A 32298-2023.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/tillsyn/A 32298-2023.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/tillsynsmail/A 32298-2023.docx")'
